$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 19702
$ws.Range("E2").Value = 1337
$ws.Range("F2").Value = 1337
$ws.Range("G2").Value = 1779
$ws.Range("H2").Value = 2563
$ws.Range("I2").Value = 2138
$ws.Range("J2").Value = 425
$ws.Range("K2").Value = 22666
$ws.Range("L2").Value = 13807
$ws.Range("M2").Value = 8859
$ws.Range("N2").Value = 4296
$ws.Range("O2").Value = 4563
$ws.Range("P2").Value = 112
$ws.Range("Q2").Value = 1150
$ws.Range("R2").Value = 1338
$ws.Range("S2").Value = -2247
$ws.Range("T2").Value = 1340
$ws.Range("U2").Value = -189
$ws.Range("V2").Value = 5382
$ws.Range("W2").Value = 6.79
$ws.Range("X2").Value = 13.01
$ws.Range("Y2").Value = 44.03
$ws.Range("Z2").Value = 12.27
$ws.Range("AA2").Value = 155.84
$ws.Range("AB2").Value = 8660.63
$ws.Range("AC2").Value = 8145
$ws.Range("AD2").Value = 1.49
$ws.Range("AE2").Value = 20702
$ws.Range("AF2").Value = 0.58
$ws.Range("AG2").Value = 150
$ws.Range("AH2").Value = 1.24
$ws.Range("AI2").Value = 1.46
$ws.Range("AJ2").Value = 22485004

# Row 3
$ws.Range("D3").Value = 30658
$ws.Range("E3").Value = 1372
$ws.Range("F3").Value = 1372
$ws.Range("G3").Value = 1278
$ws.Range("H3").Value = 942
$ws.Range("I3").Value = 406
$ws.Range("J3").Value = 536
$ws.Range("K3").Value = 25955
$ws.Range("L3").Value = 16147
$ws.Range("M3").Value = 9808
$ws.Range("N3").Value = 4724
$ws.Range("O3").Value = 5084
$ws.Range("P3").Value = 112
$ws.Range("Q3").Value = 1799
$ws.Range("R3").Value = -4274
$ws.Range("S3").Value = 2206
$ws.Range("T3").Value = 3543
$ws.Range("U3").Value = -1744
$ws.Range("V3").Value = 7868
$ws.Range("W3").Value = 4.48
$ws.Range("X3").Value = 3.07
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 3.87
$ws.Range("AA3").Value = 164.63
$ws.Range("AB3").Value = 8983.35
$ws.Range("AC3").Value = 1805
$ws.Range("AD3").Value = 7.04
$ws.Range("AE3").Value = 22766
$ws.Range("AF3").Value = 0.56
$ws.Range("AG3").Value = 150
$ws.Range("AH3").Value = 1.18
$ws.Range("AI3").Value = 7.67
$ws.Range("AJ3").Value = 22485004

# Row 4
$ws.Range("D4").Value = 31883
$ws.Range("E4").Value = 1268
$ws.Range("F4").Value = 1268
$ws.Range("G4").Value = 1702
$ws.Range("H4").Value = 1324
$ws.Range("I4").Value = 765
$ws.Range("J4").Value = 559
$ws.Range("K4").Value = 29335
$ws.Range("L4").Value = 18615
$ws.Range("M4").Value = 10720
$ws.Range("N4").Value = 5445
$ws.Range("O4").Value = 5275
$ws.Range("P4").Value = 117
$ws.Range("Q4").Value = 1266
$ws.Range("R4").Value = -2327
$ws.Range("S4").Value = 1188
$ws.Range("T4").Value = 2769
$ws.Range("U4").Value = -1503
$ws.Range("V4").Value = 9969
$ws.Range("W4").Value = 3.98
$ws.Range("X4").Value = 4.15
$ws.Range("Y4").Value = 15.04
$ws.Range("Z4").Value = 4.79
$ws.Range("AA4").Value = 173.64
$ws.Range("AB4").Value = 9313.55
$ws.Range("AC4").Value = 3394
$ws.Range("AD4").Value = 3.15
$ws.Range("AE4").Value = 25086
$ws.Range("AF4").Value = 0.43
$ws.Range("AG4").Value = 150
$ws.Range("AH4").Value = 1.4
$ws.Range("AI4").Value = 4.26
$ws.Range("AJ4").Value = 23479844

# Row 5
$ws.Range("D5").Value = 29132
$ws.Range("E5").Value = -258
$ws.Range("F5").Value = -258
$ws.Range("G5").Value = -913
$ws.Range("H5").Value = -222
$ws.Range("I5").Value = -146
$ws.Range("J5").Value = -75
$ws.Range("K5").Value = 27677
$ws.Range("L5").Value = 17677
$ws.Range("M5").Value = 10000
$ws.Range("N5").Value = 5070
$ws.Range("O5").Value = 4930
$ws.Range("P5").Value = 117
$ws.Range("Q5").Value = 1910
$ws.Range("R5").Value = -1352
$ws.Range("S5").Value = 253
$ws.Range("T5").Value = 2052
$ws.Range("U5").Value = -142
$ws.Range("V5").Value = 9457
$ws.Range("W5").Value = -0.88
$ws.Range("X5").Value = -0.76
$ws.Range("Y5").Value = -2.78
$ws.Range("Z5").Value = -0.78
$ws.Range("AA5").Value = 176.76
$ws.Range("AB5").Value = 9156.81
$ws.Range("AC5").Value = -623
$ws.Range("AD5").Value = -11.82
$ws.Range("AE5").Value = 23358
$ws.Range("AF5").Value = 0.32
$ws.Range("AG5").Value = 100
$ws.Range("AH5").Value = 1.36
$ws.Range("AI5").Value = -14.84
$ws.Range("AJ5").Value = 23479844

# Row 6
$ws.Range("D6").Value = 28525
$ws.Range("E6").Value = 70
$ws.Range("F6").Value = 70
$ws.Range("G6").Value = -164
$ws.Range("H6").Value = -432
$ws.Range("I6").Value = -178
$ws.Range("K6").Value = 26963
$ws.Range("L6").Value = 17460
$ws.Range("M6").Value = 9503
$ws.Range("N6").Value = 4789
$ws.Range("P6").Value = 117
$ws.Range("Q6").Value = 536
$ws.Range("R6").Value = -2720
$ws.Range("S6").Value = 393
$ws.Range("T6").Value = 2522
$ws.Range("U6").Value = -1986
$ws.Range("V6").Value = 9369
$ws.Range("W6").Value = 0.25
$ws.Range("X6").Value = -1.51
$ws.Range("Y6").Value = -3.61
$ws.Range("Z6").Value = -1.58
$ws.Range("AA6").Value = 183.74
$ws.Range("AB6").Value = 8929.68
$ws.Range("AC6").Value = -757
$ws.Range("AD6").Value = -4.79
$ws.Range("AE6").Value = 22063
$ws.Range("AF6").Value = 0.16
$ws.Range("AG6").Value = 50
$ws.Range("AH6").Value = 1.38
$ws.Range("AI6").Value = -6.1
$ws.Range("AJ6").Value = 23479844

# Clear row 7
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Clear row 8
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Clear row 9
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
